$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Annual")
$ws2.Columns.Item(60).ColumnWidth = 9.8585
